$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, [string]$innerXml) {
    $r = $paragraph.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Delete()
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# 1) ORDERDETAILS -> ORDERDETAIL (global rename; affects the two simple paragraphs
#    as well as the paragraph that gets fully rebuilt below)
$d.Content.Find.Execute("ORDERDETAILS", $true, $false, $false, $false, $false, $true, 1, $false, "ORDERDETAIL", 2) | Out-Null

# 2) "ORDERDETAIL (ticket) -> ManyToOne -> TICKET (orderDetails)" paragraph
#    becomes "ORDERDETAIL (product) -> ManyToOne -> PRODUCT (orderDetails)"
$p10 = $d.Paragraphs.Item(10)
Set-ParagraphXml $p10 '<w:r><w:t>ORDERDETAIL (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> -&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ManyToOne</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-&gt; </w:t></w:r><w:r><w:t>PRODUCT</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>orderDetails</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r>'

# 3) "TICKETPRICING (ticket) -> ManyToOne ->TICKET (ticketPricings)" paragraph
#    becomes "PRODUCTPRICING (product) -> ManyToOne -> PRODUCT (productPricings)"
#    with the closing parenthesis highlighted in yellow
$p12 = $d.Paragraphs.Item(12)
Set-ParagraphXml $p12 '<w:r><w:t>PRODUCT</w:t></w:r><w:r><w:t xml:space="preserve">PRICING </w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) </w:t></w:r><w:r><w:t xml:space="preserve">-&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ManyToOne</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>-&gt;</w:t></w:r><w:r><w:t>PRODUCT</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>product</w:t></w:r><w:r><w:t>Pricing</w:t></w:r><w:r><w:t>s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>)</w:t></w:r>'

# 4) "TICKETPRICING (exhibition) -> ManyToOne ->EXHIBITION (ticketPricings)" paragraph
#    becomes "PRODUCTPRICING (exhibition) -> ManyToOne -> EXHIBITION (product Pricings)"
$p13 = $d.Paragraphs.Item(13)
Set-ParagraphXml $p13 '<w:r><w:t xml:space="preserve">PRODUCTPRICING </w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:t>exhibition</w:t></w:r><w:r><w:t xml:space="preserve">) </w:t></w:r><w:r><w:t xml:space="preserve">-&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ManyToOne</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>-&gt;</w:t></w:r><w:r><w:t>EXHIBITION</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pri</w:t></w:r><w:r><w:t>cing</w:t></w:r><w:r><w:t>s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r>'

# 5) After the "COMMENT (exhibition) -> ManyToOne -> EXHIBITION (comments)" paragraph,
#    add a blank paragraph followed by a new
#    "PRODUCT (types) -> OneToMany -> TYPE (product)" paragraph
$pComment = $d.Paragraphs.Item(19)
$pComment.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs.Item(20)
$pBlank.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs.Item(21)
Set-ParagraphXml $pNew '<w:r><w:t xml:space="preserve">PRODUCT (types) -&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>OneToMany</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -&gt; TYPE (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r>'
